# Refresh the cryptos table (Price = column D, Volume(1h) = column E) with the
# latest scrape. All figures are plain text in the sheet (not numbers/formulas),
# so values are written as strings. A handful of Price cells now read as plain
# decimals (e.g. "312.30"); Excel would otherwise auto-convert those to numbers
# (dropping the trailing zero), so we briefly force Text format for the write and
# then restore the normal style so no formatting change is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.334.16"
$ws.Range("E2").Value = "  -3.01%  "
$ws.Range("D3").Value = "2.466.71"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.30%  "
$ws.Range("E7").Value = "  -3.03%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0783"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.75%  "
$ws.Range("D14").Value = "2.846.03"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "2.467.07"
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.27%  "
$ws.Range("D18").Value = "41.292.06"
$ws.Range("E18").Value = "  -3.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.45%  "
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("E21").Value = "  -5.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.95%  "
$ws.Range("E25").Value = "  -5.10%  "
$ws.Range("E27").Value = "  -5.95%  "
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("E33").Value = "  -6.04%  "
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0756"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("E36").Value = "  -5.50%  "
$ws.Range("E37").Value = "  -6.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.06%  "
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("E40").Value = "  -7.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.33%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "1.965.09"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").Value = "  -4.51%  "
$ws.Range("E46").Value = "  -6.39%  "
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "70.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.42%  "
$ws.Range("E51").Value = "  -5.31%  "
